$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "level": rework the arrow/bullet/special per-row spawn markers
# (levels raised up to 50 -> updated layout for rows 1-19)
# ---------------------------------------------------------------------
$level = $wb.Worksheets.Item("level")

# Clear cells that no longer hold a marker
$clearAddrs = @("B1","A2","A4","B4","C4","D5","A8","B8","C9","B12","D12","A15","D15","B17","D18","B19")
foreach ($addr in $clearAddrs) {
    $level.Range($addr).ClearContents()
}

# Set cells with their new marker values
$level.Range("A1").Value = 1
$level.Range("B2").Value = 2
$level.Range("D2").Value = 1
$level.Range("C3").Value = 8
$level.Range("B5").Value = 4
$level.Range("A7").Value = 8
$level.Range("C7").Value = 9
$level.Range("D8").Value = 6
$level.Range("D10").Value = 9
$level.Range("B11").Value = 8
$level.Range("C12").Value = 7

# Apply the new "vertical center" style to every cell that is actually in use
$styledAddrs = @(
    "A1","E1",
    "B2","D2","E2",
    "C3","E3",
    "E4",
    "B5","E5",
    "E6",
    "A7","C7","E7",
    "D8","E8",
    "E9",
    "D10","E10",
    "B11","E11",
    "C12","E12"
)
for ($r = 13; $r -le 48; $r++) {
    $styledAddrs += "E$r"
}
foreach ($addr in $styledAddrs) {
    $level.Range($addr).VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# Sheet "enemies": drop the stray group notes, extend the table to row 10
# ---------------------------------------------------------------------
$enemies = $wb.Worksheets.Item("enemies")
$enemies.Range("G6").ClearContents()
$enemies.Range("G7").ClearContents()

$enemies.Range("A8").Value = 7
$enemies.Range("B8").Value = 4
$enemies.Range("C8").Value = 0
$enemies.Range("D8").Value = 0
$enemies.Range("E8").Value = 1
$enemies.Range("F8").Value = "EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow"

$enemies.Range("A9").Value = 8
$enemies.Range("B9").Value = 0
$enemies.Range("C9").Value = 2
$enemies.Range("D9").Value = 0
$enemies.Range("E9").Value = 1
$enemies.Range("F9").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"

$enemies.Range("A10").Value = 9
$enemies.Range("B10").Value = 0
$enemies.Range("C10").Value = 0
$enemies.Range("D10").Value = 2
$enemies.Range("E10").Value = 1
$enemies.Range("F10").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"

foreach ($addr in @("A8","A9","A10")) {
    $enemies.Range($addr).HorizontalAlignment = -4108
    $enemies.Range($addr).VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")
$meta.Range("B39").Select() | Out-Null

$enemies.Range("F8").Select() | Out-Null

$misc = $wb.Worksheets.Item("misc")
$misc.Range("A49").Select() | Out-Null

$level.Activate() | Out-Null
$level.Range("E6").Select() | Out-Null
